$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.139.76"
$ws.Range("E2").Value = "  -4.33%  "
$ws.Range("D3").Value = "1.655.10"
$ws.Range("E3").Value = "  -3.26%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "215.46"
$ws.Range("E5").Value = "  -3.98%  "
$ws.Range("E6").Value = "  -3.80%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "0.2580"
$ws.Range("D9").Value = "0.06411"
$ws.Range("E9").Value = "  -4.49%  "
$ws.Range("D10").Value = "19.89"
$ws.Range("E10").Value = "  -4.62%  "
$ws.Range("D11").Value = "0.07796"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").Value = "1.658.18"
$ws.Range("E12").Value = "  -3.09%  "
$ws.Range("D14").Value = "1.882.38"
$ws.Range("E14").Value = "  -3.28%  "
$ws.Range("D15").Value = "0.5512"
$ws.Range("E15").Value = "  -5.29%  "
$ws.Range("D16").Value = "0.0₅8006"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").Value = "63.98"
$ws.Range("E17").Value = "  -6.20%  "
$ws.Range("D18").Value = "26.155.28"
$ws.Range("E18").Value = "  -4.32%  "
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "209.08"
$ws.Range("E20").Value = "  -7.56%  "
$ws.Range("D21").Value = "4.408"
$ws.Range("E21").Value = "  -4.69%  "
$ws.Range("D22").Value = "10.06"
$ws.Range("E22").Value = "  -3.22%  "
$ws.Range("D23").Value = "6.039"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").Value = "143.76"
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("D27").Value = "0.1179"
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("D28").Value = "6.974"
$ws.Range("E28").Value = "  -3.62%  "
$ws.Range("E29").Value = "  -3.25%  "
$ws.Range("D30").Value = "0.05092"
$ws.Range("E30").Value = "  -5.00%  "
$ws.Range("E31").Value = "  -3.81%  "
$ws.Range("D32").Value = "3.344"
$ws.Range("E32").Value = "  -3.88%  "
$ws.Range("D33").Value = "3.218"
$ws.Range("E33").Value = "  -6.18%  "
$ws.Range("D34").Value = "1.566"
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("D35").Value = "2.753"
$ws.Range("E35").Value = "  -4.09%  "
$ws.Range("D36").Value = "0.9284"
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("D37").Value = "2.365"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").Value = "1.164.39"
$ws.Range("E38").Value = "  +7.65%  "
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("E40").Value = "  -2.93%  "
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "5.652"
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("D45").Value = "100.44"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").Value = "1.792.37"
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").Value = "0.4551"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("E49").Value = "  -3.44%  "
$ws.Range("D50").Value = "1.007"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("D51").Value = "7.834"
$ws.Range("E51").Value = "  -2.97%  "
